# daily auto push: 2026-02-10 14:24 UTC
#
# A new reading was recorded for 2026/02/10 (weekday 火) at hour 21.
# It belongs right after the existing 2026/02/10 rows (currently at
# sheet rows 793-794) and before the 2026/12/29 block that currently
# starts at row 795, so insert a new row there and push everything
# else down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 795; this shifts rows
# 795-836 down to 796-837 and keeps rows 1-794 untouched.
$ws.Rows.Item(795).Insert()

# Fill in the newly inserted row 795 with the new record.
# Column A stores a date formatted as plain text (e.g. "2026/02/10"),
# matching every other date cell in the sheet (they are literal
# strings, not Excel date serials). Temporarily force a text number
# format before assigning the value so the COM layer doesn't
# auto-convert the "yyyy/mm/dd" string into a real date, then clear
# the temporary formatting so the cell is left with no explicit style
# (matching the unstyled cells around it).
$ws.Range("A795").NumberFormat = "@"
$ws.Range("A795").Value = "2026/02/10"
$ws.Range("A795").ClearFormats()

$ws.Range("B795").Value = "火"
$ws.Range("C795").Value = 21
$ws.Range("D795").Value = 201
